# Update math problems in the table (three-digit divided by one-digit).
# Each Find/Replace is scoped to its specific table cell to target the exact
# occurrence; operations are ordered so that a replacement never creates text
# that would be mistakenly matched by a not-yet-executed replacement elsewhere
# in the document (Find searches the whole document for the first/leftmost
# match, even when invoked on a single cell's Range).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.Find.Execute("205÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "734÷9=", 2) | Out-Null

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.Find.Execute("746÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "986÷4=", 2) | Out-Null

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.Find.Execute("112÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "434÷5=", 2) | Out-Null

$cell = $t.Cell(1, 4)
$rng = $cell.Range
$rng.Find.Execute("548÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "155÷4=", 2) | Out-Null

$cell = $t.Cell(1, 5)
$rng = $cell.Range
$rng.Find.Execute("445÷5=", $false, $false, $false, $false, $false, $true, 0, $false, "407÷8=", 2) | Out-Null

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.Find.Execute("182÷8=", $false, $false, $false, $false, $false, $true, 0, $false, "288÷8=", 2) | Out-Null

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.Find.Execute("837÷5=", $false, $false, $false, $false, $false, $true, 0, $false, "360÷8=", 2) | Out-Null

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.Find.Execute("945÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "935÷4=", 2) | Out-Null

$cell = $t.Cell(5, 4)
$rng = $cell.Range
$rng.Find.Execute("466÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "749÷3=", 2) | Out-Null

$cell = $t.Cell(5, 5)
$rng = $cell.Range
$rng.Find.Execute("667÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "464÷2=", 2) | Out-Null

$cell = $t.Cell(9, 1)
$rng = $cell.Range
$rng.Find.Execute("908÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "781÷9=", 2) | Out-Null

$cell = $t.Cell(13, 5)
$rng = $cell.Range
$rng.Find.Execute("224÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "439÷2=", 2) | Out-Null

$cell = $t.Cell(9, 2)
$rng = $cell.Range
$rng.Find.Execute("528÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "224÷6=", 2) | Out-Null

$cell = $t.Cell(9, 3)
$rng = $cell.Range
$rng.Find.Execute("864÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "727÷3=", 2) | Out-Null

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.Find.Execute("620÷6=", $false, $false, $false, $false, $false, $true, 0, $false, "852÷6=", 2) | Out-Null

$cell = $t.Cell(9, 5)
$rng = $cell.Range
$rng.Find.Execute("499÷8=", $false, $false, $false, $false, $false, $true, 0, $false, "246÷8=", 2) | Out-Null

$cell = $t.Cell(13, 1)
$rng = $cell.Range
$rng.Find.Execute("330÷8=", $false, $false, $false, $false, $false, $true, 0, $false, "711÷2=", 2) | Out-Null

$cell = $t.Cell(13, 2)
$rng = $cell.Range
$rng.Find.Execute("436÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "526÷5=", 2) | Out-Null

$cell = $t.Cell(13, 3)
$rng = $cell.Range
$rng.Find.Execute("705÷2=", $false, $false, $false, $false, $false, $true, 0, $false, "906÷6=", 2) | Out-Null

$cell = $t.Cell(13, 4)
$rng = $cell.Range
$rng.Find.Execute("268÷9=", $false, $false, $false, $false, $false, $true, 0, $false, "818÷8=", 2) | Out-Null

$cell = $t.Cell(17, 1)
$rng = $cell.Range
$rng.Find.Execute("952÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "240÷6=", 2) | Out-Null

$cell = $t.Cell(17, 2)
$rng = $cell.Range
$rng.Find.Execute("847÷4=", $false, $false, $false, $false, $false, $true, 0, $false, "598÷5=", 2) | Out-Null

$cell = $t.Cell(17, 3)
$rng = $cell.Range
$rng.Find.Execute("303÷3=", $false, $false, $false, $false, $false, $true, 0, $false, "792÷2=", 2) | Out-Null

$cell = $t.Cell(17, 4)
$rng = $cell.Range
$rng.Find.Execute("806÷7=", $false, $false, $false, $false, $false, $true, 0, $false, "421÷7=", 2) | Out-Null

$cell = $t.Cell(17, 5)
$rng = $cell.Range
$rng.Find.Execute("319÷9=", $false, $false, $false, $false, $false, $true, 0, $false, "325÷5=", 2) | Out-Null
